$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.360.59"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.52%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.843.70"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.84%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9987"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.37%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.19"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.76%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6297"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.73%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.0000"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.39%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07438"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.04%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2906"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.81%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.83"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.83%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07742"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.39%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.847.88"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.55%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.986"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.27%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6795"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.98%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001026"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.66%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "82.07"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.73%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.276"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.85%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.338.62"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.65%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "229.48"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.49%  "

# Row 20
$ws.Range("E20").Value = "  -0.73%  "

# Row 21
$ws.Range("E21").Value = "  -0.38%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.443"
$ws.Range("D22").ClearFormats()

# Row 23
$ws.Range("E23").Value = "  -0.29%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "158.26"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.91%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.494"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.15%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1356"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.87%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.45"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.75%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06492"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +13.63%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.448"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.76%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.489"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.42%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.072"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.13%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.062"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.33%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.840"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.23%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.140"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.89%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6959"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.33%  "

# Row 36
$ws.Range("E36").Value = "  -0.95%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01856"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.05%  "

# Row 38
$ws.Range("E38").Value = "  +1.27%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.241.59"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.39%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.804"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +4.32%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9323"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.51%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9996"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.44%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.998.18"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.02%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.92"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.82%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "65.64"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.09%  "

# Row 46
$ws.Range("E46").Value = "  +2.51%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.053"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.43%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.713"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.92%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1152"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.67%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.998"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.74%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3894"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.31%  "
